$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new product row above the current row 7 ("BOBAI SUNSCREEN ...").
# This pushes the existing product rows (and the totals/footer rows below
# them) down by one, exactly like Excel's normal "Insert Row" behaviour.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Insert()

# Copy the formatting (styles, number formats, fonts, fills, borders) of the
# row that used to be row 7 and is now row 8 into the freshly inserted row 7
# so the new row looks identical to the other product rows.
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height used by the other product rows.
$ws.Rows.Item(7).RowHeight = 25.5

# Recreate the merged cells for the new row (mirrors A8:B8, C8:G8, H8:K8,
# L8:M8, N8:O8 on the sibling rows).
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# ---------------------------------------------------------------------------
# Populate the new row with the "ATENO 50MG 20 F.C.TAB." entry.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "ATENO 50MG 20 F.C.TAB."

# The ratio / quantity / price columns are stored as *text* in this sheet
# (even though they look numeric), so force text formatting before writing
# them to avoid Excel auto-converting them to numbers.
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "2:1"

$fmtL = $ws.Range("L8").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $fmtL

$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "12.00"

$fmtP = $ws.Range("P8").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "6.0000"
$ws.Range("P7").NumberFormat = $fmtP

$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "0:1"

# ---------------------------------------------------------------------------
# Update the running total of transactions shown at the bottom of the sheet
# (now on row 12 after the insert) from 581 to 587.
# ---------------------------------------------------------------------------
$ws.Range("P12").Value = 587
